
# ============================================================================
# CREACION DE VENTANA DE ERROR, SELECT CON BUSCADOR EN LA VENTANA REGISTRO
#
# - Clears the stray duplicated "Carga al servidor de producción" entry
#   that had been left in B12 (row stays part of the A6:A12 merge, just the
#   text is gone now).
# - Updates the "login" validation note (row 13) to its final wording with a
#   trailing period.
# - Adds two new activity rows at the bottom of the log:
#     14: "*Implementación de un Select con input en el formulario de
#          registro."  (marked with an X under ACTUALIZACION)
#     15: "*Maquetación de  ventana de error."
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: the duplicate task line is removed, row/merge stays in place ---
$ws.Range("B12").ClearContents()

# --- Make room for the two new activity rows under row 13 -------------------
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(15).Insert()

# Drop the formatting that Insert() copied down into the blank cells of the
# two new rows - in the final sheet these rows only carry the cells that
# actually hold data (B14/D14 and B15).
$ws.Range("A14:E15").Clear()

# --- Fill in the new activity rows ------------------------------------------
$ws.Range("B15").Value = "*Maquetación de  ventana de error."
$ws.Range("B14").Value = "*Implementación de un Select con input en el formulario de registro."
$ws.Range("D14").Value = "X"

# --- Final wording for the login validation task on row 13 ------------------
$ws.Range("B13").Value = "* Validación de los campos correo y clave en el login."

# --- Match the workbook's last on-screen selection ---------------------------
$ws.Range("B21").Select() | Out-Null
